$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 6.173199999999998
$ws.Range("D4").Value = -8.215299999999996
$ws.Range("E4").Value = 13.5203

$ws.Range("D5").Value = -8.054399999999996

$ws.Range("B7").Value = 6.852599999999995

$ws.Range("D8").Value = -8.350599999999993

$ws.Range("E9").Value = 13.90110000000001

$ws.Range("B16").Value = 8.576800000000011
$ws.Range("D16").Value = -8.277699999999998

$ws.Range("E18").Value = 13.1154
